$d = $word.ActiveDocument

# 1) Reorder skills list for 2007 - 2011 entry: move "Paypal Payflow Pro" after "GIT"
$d.Content.Find.Execute(
    "2007 - 2011 — Python, Django, Javascript, Backbone.js, MySQL, SOLR, RabbitMQ, Amazon EC2, Paypal Payflow Pro, GIT",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2007 - 2011 — Python, Django, Javascript, Backbone.js, MySQL, SOLR, RabbitMQ, Amazon EC2, GIT, Paypal Payflow Pro",
    2
)

# 2) Reorder skills list for 2007 entry: move "Paypal Payflow Pro" after "SVN"
$d.Content.Find.Execute(
    "2007 — Python, Django, Javascript, Backbone.js, MySQL, Xapian, Amazon EC2, Paypal Payflow Pro, GIT, SVN",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2007 — Python, Django, Javascript, Backbone.js, MySQL, Xapian, Amazon EC2, GIT, SVN, Paypal Payflow Pro",
    2
)

# 3) Reorder skills list for 1998 - 2004 entry: move "DB2, Oracle, Microsoft IIS" after "Javascript"
$d.Content.Find.Execute(
    "1998 - 2004 — Java, JSP, PHP, ASP, Visual Basic, DB2, Oracle, Microsoft IIS, HTML, CSS, Javascript",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "1998 - 2004 — Java, JSP, PHP, ASP, Visual Basic, HTML, CSS, Javascript, DB2, Oracle, Microsoft IIS",
    2
)
